# Update "countries & provincias Spain" workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 14:22"

# Swap the city names for rows 24 and 25 (Sevilla <-> Gipuzkoa/Guipuzcoa)
$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("A25").Value = "Sevilla"

# Row 7 (Bizkaia/Vizcaya) updated figures
$ws.Range("B7").Value = 7045
$ws.Range("C7").Value = 7124
$ws.Range("D7").Value = 4423
$ws.Range("E7").Value = 551

# Row 16 (Araba/Alava) updated figures
$ws.Range("B16").Value = 3241
$ws.Range("C16").Value = 7124
$ws.Range("D16").Value = 4423
$ws.Range("E16").Value = 318

# Row 24 (now Gipuzkoa/Guipuzcoa) updated figures
$ws.Range("B24").Value = 2342
$ws.Range("C24").Value = 7124
$ws.Range("D24").Value = 4423
$ws.Range("E24").Value = 212

# Row 25 (now Sevilla) updated figures
$ws.Range("B25").Value = 2329
$ws.Range("C25").Value = 459
$ws.Range("D25").Value = 1658
$ws.Range("E25").Value = 212
